{"js": "const replacements = [\n  [\"953\u00f75=190, 3\", \"235\u00f79=26, 1\"],\n  [\"849\u00f79=94, 3\", \"780\u00f76=130, 0\"],\n  [\"641\u00f79=71, 2\", \"469\u00f79=52, 1\"],\n  [\"458\u00f77=65, 3\", \"828\u00f73=276, 0\"],\n  [\"582\u00f74=145, 2\", \"924\u00f78=115, 4\"],\n  [\"110\u00f76=18, 2\", \"539\u00f78=67, 3\"],\n  [\"214\u00f78=26, 6\", \"823\u00f76=137, 1\"],\n  [\"611\u00f73=203, 2\", \"896\u00f79=99, 5\"],\n  [\"751\u00f74=187, 3\", \"105\u00f75=21, 0\"],\n  [\"180\u00f76=30, 0\", \"214\u00f76=35, 4\"],\n  [\"618\u00f76=103, 0\", \"312\u00f78=39, 0\"],\n  [\"152\u00f75=30, 2\", \"781\u00f72=390, 1\"],\n  [\"340\u00f72=170, 0\", \"879\u00f78=109, 7\"],\n  [\"356\u00f76=59, 2\", \"123\u00f73=41, 0\"],\n  [\"185\u00f74=46, 1\", \"291\u00f76=48, 3\"],\n  [\"142\u00f73=47, 1\", \"120\u00f76=20, 0\"],\n  [\"719\u00f79=79, 8\", \"153\u00f73=51, 0\"],\n  [\"901\u00f72=450, 1\", \"463\u00f73=154, 1\"],\n  [\"123\u00f77=17, 4\", \"670\u00f75=134, 0\"],\n  [\"988\u00f73=329, 1\", \"341\u00f77=48, 5\"],\n  [\"416\u00f76=69, 2\", \"899\u00f75=179, 4\"],\n  [\"904\u00f77=129, 1\", \"887\u00f72=443, 1\"],\n  [\"163\u00f72=81, 1\", \"796\u00f75=159, 1\"],\n  [\"512\u00f77=73, 1\", \"846\u00f74=211, 2\"],\n  [\"219\u00f79=24, 3\", \"299\u00f76=49, 5\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"953\u00f75=190, 3\", \"235\u00f79=26, 1\"),\n    @(\"849\u00f79=94, 3\", \"780\u00f76=130, 0\"),\n    @(\"641\u00f79=71, 2\", \"469\u00f79=52, 1\"),\n    @(\"458\u00f77=65, 3\", \"828\u00f73=276, 0\"),\n    @(\"582\u00f74=145, 2\", \"924\u00f78=115, 4\"),\n    @(\"110\u00f76=18, 2\", \"539\u00f78=67, 3\"),\n    @(\"214\u00f78=26, 6\", \"823\u00f76=137, 1\"),\n    @(\"611\u00f73=203, 2\", \"896\u00f79=99, 5\"),\n    @(\"751\u00f74=187, 3\", \"105\u00f75=21, 0\"),\n    @(\"180\u00f76=30, 0\", \"214\u00f76=35, 4\"),\n    @(\"618\u00f76=103, 0\", \"312\u00f78=39, 0\"),\n    @(\"152\u00f75=30, 2\", \"781\u00f72=390, 1\"),\n    @(\"340\u00f72=170, 0\", \"879\u00f78=109, 7\"),\n    @(\"356\u00f76=59, 2\", \"123\u00f73=41, 0\"),\n    @(\"185\u00f74=46, 1\", \"291\u00f76=48, 3\"),\n    @(\"142\u00f73=47, 1\", \"120\u00f76=20, 0\"),\n    @(\"719\u00f79=79, 8\", \"153\u00f73=51, 0\"),\n    @(\"901\u00f72=450, 1\", \"463\u00f73=154, 1\"),\n    @(\"123\u00f77=17, 4\", \"670\u00f75=134, 0\"),\n    @(\"988\u00f73=329, 1\", \"341\u00f77=48, 5\"),\n    @(\"416\u00f76=69, 2\", \"899\u00f75=179, 4\"),\n    @(\"904\u00f77=129, 1\", \"887\u00f72=443, 1\"),\n    @(\"163\u00f72=81, 1\", \"796\u00f75=159, 1\"),\n    @(\"512\u00f77=73, 1\", \"846\u00f74=211, 2\"),\n    @(\"219\u00f79=24, 3\", \"299\u00f76=49, 5\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
